$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 972-1061 (date serial, age group, covid_deaths)
# This reflects a refresh of the auto-generated dataset: historical counts
# were revised and four new rows (1058-1061, date 2020-12-09) were appended.
$rows = @(
    @(972, 44156, "60-69", 9),
    @(973, 44156, "70-79", 12),
    @(974, 44156, "80+", 35),
    @(975, 44157, "30-39", 1),
    @(976, 44157, "50-59", 2),
    @(977, 44157, "60-69", 8),
    @(978, 44157, "70-79", 17),
    @(979, 44157, "80+", 36),
    @(980, 44158, "60-69", 13),
    @(981, 44158, "70-79", 24),
    @(982, 44158, "80+", 35),
    @(983, 44159, "50-59", 1),
    @(984, 44159, "60-69", 9),
    @(985, 44159, "70-79", 17),
    @(986, 44159, "80+", 34),
    @(987, 44160, "30-39", 1),
    @(988, 44160, "40-49", 1),
    @(989, 44160, "50-59", 1),
    @(990, 44160, "60-69", 6),
    @(991, 44160, "70-79", 15),
    @(992, 44160, "80+", 31),
    @(993, 44161, "40-49", 3),
    @(994, 44161, "50-59", 4),
    @(995, 44161, "60-69", 9),
    @(996, 44161, "70-79", 15),
    @(997, 44161, "80+", 38),
    @(998, 44162, "50-59", 4),
    @(999, 44162, "60-69", 8),
    @(1000, 44162, "70-79", 31),
    @(1001, 44162, "80+", 28),
    @(1002, 44163, "30-39", 2),
    @(1003, 44163, "40-49", 1),
    @(1004, 44163, "50-59", 3),
    @(1005, 44163, "60-69", 13),
    @(1006, 44163, "70-79", 12),
    @(1007, 44163, "80+", 29),
    @(1008, 44164, "50-59", 2),
    @(1009, 44164, "60-69", 8),
    @(1010, 44164, "70-79", 18),
    @(1011, 44164, "80+", 36),
    @(1012, 44165, "50-59", 2),
    @(1013, 44165, "60-69", 6),
    @(1014, 44165, "70-79", 20),
    @(1015, 44165, "80+", 37),
    @(1016, 44166, "40-49", 1),
    @(1017, 44166, "50-59", 2),
    @(1018, 44166, "60-69", 8),
    @(1019, 44166, "70-79", 20),
    @(1020, 44166, "80+", 33),
    @(1021, 44167, "40-49", 2),
    @(1022, 44167, "50-59", 5),
    @(1023, 44167, "60-69", 7),
    @(1024, 44167, "70-79", 16),
    @(1025, 44167, "80+", 45),
    @(1026, 44168, "40-49", 1),
    @(1027, 44168, "50-59", 1),
    @(1028, 44168, "60-69", 7),
    @(1029, 44168, "70-79", 20),
    @(1030, 44168, "80+", 36),
    @(1031, 44169, "30-39", 1),
    @(1032, 44169, "40-49", 3),
    @(1033, 44169, "50-59", 3),
    @(1034, 44169, "60-69", 12),
    @(1035, 44169, "70-79", 17),
    @(1036, 44169, "80+", 36),
    @(1037, 44170, "50-59", 3),
    @(1038, 44170, "60-69", 12),
    @(1039, 44170, "70-79", 16),
    @(1040, 44170, "80+", 37),
    @(1041, 44171, "20-29", 1),
    @(1042, 44171, "50-59", 2),
    @(1043, 44171, "60-69", 10),
    @(1044, 44171, "70-79", 18),
    @(1045, 44171, "80+", 23),
    @(1046, 44172, "30-39", 1),
    @(1047, 44172, "40-49", 1),
    @(1048, 44172, "50-59", 3),
    @(1049, 44172, "60-69", 14),
    @(1050, 44172, "70-79", 19),
    @(1051, 44172, "80+", 34),
    @(1052, 44173, "0-19", 1),
    @(1053, 44173, "20-29", 1),
    @(1054, 44173, "50-59", 2),
    @(1055, 44173, "60-69", 11),
    @(1056, 44173, "70-79", 12),
    @(1057, 44173, "80+", 42),
    @(1058, 44174, "50-59", 1),
    @(1059, 44174, "60-69", 2),
    @(1060, 44174, "70-79", 8),
    @(1061, 44174, "80+", 23)
)

$dateFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateCell = $ws.Cells.Item($rowNum, 1)
    $dateCell.Value = $r[1]
    $dateCell.NumberFormat = $dateFormat
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}
